$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells in H1/I1, styled like the rest of row 1 (bold, centered)
$ws.Range("H1").Value = "Threat_total"
$ws.Range("I1").Value = "Threat_prc"
$ws.Range("H1:I1").Font.Bold = $true
$ws.Range("H1:I1").HorizontalAlignment = -4108  # xlCenter

# Data rows:
#   Threat_total = Vulnerable + Endangered + Critically Endangered
#   Threat_prc   = round(Threat_total / Total * 100)
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = 75

$ws.Range("H3").Value = 6
$ws.Range("I3").Value = 50

$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 75

$ws.Range("H5").Value = 11
$ws.Range("I5").Value = 65

$ws.Range("H6").Value = 6
$ws.Range("I6").Value = 75
